# Append the new "Dani" reference row (row 53) to Sheet1 and turn on the
# header AutoFilter, matching the "Updating content to new structure" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the new data row -------------------------------------------------
$ws.Range("A53").Value = "Dani"
$ws.Range("B53").Value = "dani.jpeg"
$ws.Range("C53").Value = "dani_0.jpeg"
$ws.Range("D53").Value = "dani_1.jpeg"
$ws.Range("E53").Value = "dani_2.jpeg"
$ws.Range("F53").Value = "dani_3.jpeg"

# --- 2. Turn on the AutoFilter for the header row ---------------------------
$ws.Range("A1:F1").AutoFilter()

# --- 3. Register the (hidden) _FilterDatabase defined name that Excel writes
#        whenever AutoFilter is turned on for a sheet ------------------------
$fdName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$F`$1")
$fdName.Visible = $false

# --- 4. Leave the new row selected, like the final user action in Excel -----
$ws.Range("A53:F53").Select()
